{"js": "const replacements = [\n  [\"2024-08-07 Wednesday\", \"2024-08-08 Thursday\"],\n  [\"71\u00d767=4757\", \"69\u00d787=6003\"],\n  [\"91\u00d785=7735\", \"78\u00d719=1482\"],\n  [\"78\u00d747=3666\", \"82\u00d782=6724\"],\n  [\"95\u00d768=6460\", \"14\u00d714=196\"],\n  [\"88\u00d754=4752\", \"75\u00d750=3750\"],\n  [\"65\u00d746=2990\", \"71\u00d770=4970\"],\n  [\"53\u00d769=3657\", \"21\u00d756=1176\"],\n  [\"64\u00d745=2880\", \"26\u00d777=2002\"],\n  [\"34\u00d746=1564\", \"64\u00d717=1088\"],\n  [\"62\u00d786=5332\", \"59\u00d721=1239\"],\n  [\"11\u00d788=968\", \"57\u00d786=4902\"],\n  [\"11\u00d747=517\", \"32\u00d743=1376\"],\n  [\"41\u00d732=1312\", \"19\u00d779=1501\"],\n  [\"61\u00d735=2135\", \"57\u00d785=4845\"],\n  [\"95\u00d732=3040\", \"54\u00d766=3564\"],\n  [\"42\u00d799=4158\", \"11\u00d791=1001\"],\n  [\"86\u00d780=6880\", \"48\u00d718=864\"],\n  [\"88\u00d746=4048\", \"23\u00d722=506\"],\n  [\"11\u00d753=583\", \"82\u00d793=7626\"],\n  [\"94\u00d743=4042\", \"20\u00d721=420\"],\n  [\"38\u00d715=570\", \"84\u00d773=6132\"],\n  [\"71\u00d743=3053\", \"45\u00d747=2115\"],\n  [\"46\u00d779=3634\", \"76\u00d769=5244\"],\n  [\"33\u00d783=2739\", \"33\u00d769=2277\"],\n  [\"55\u00d769=3795\", \"13\u00d711=143\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (const range of searchResults.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2024-08-07 Wednesday\"; New = \"2024-08-08 Thursday\"},\n    @{Old = \"71\u00d767=4757\"; New = \"69\u00d787=6003\"},\n    @{Old = \"91\u00d785=7735\"; New = \"78\u00d719=1482\"},\n    @{Old = \"78\u00d747=3666\"; New = \"82\u00d782=6724\"},\n    @{Old = \"95\u00d768=6460\"; New = \"14\u00d714=196\"},\n    @{Old = \"88\u00d754=4752\"; New = \"75\u00d750=3750\"},\n    @{Old = \"65\u00d746=2990\"; New = \"71\u00d770=4970\"},\n    @{Old = \"53\u00d769=3657\"; New = \"21\u00d756=1176\"},\n    @{Old = \"64\u00d745=2880\"; New = \"26\u00d777=2002\"},\n    @{Old = \"34\u00d746=1564\"; New = \"64\u00d717=1088\"},\n    @{Old = \"62\u00d786=5332\"; New = \"59\u00d721=1239\"},\n    @{Old = \"11\u00d788=968\"; New = \"57\u00d786=4902\"},\n    @{Old = \"11\u00d747=517\"; New = \"32\u00d743=1376\"},\n    @{Old = \"41\u00d732=1312\"; New = \"19\u00d779=1501\"},\n    @{Old = \"61\u00d735=2135\"; New = \"57\u00d785=4845\"},\n    @{Old = \"95\u00d732=3040\"; New = \"54\u00d766=3564\"},\n    @{Old = \"42\u00d799=4158\"; New = \"11\u00d791=1001\"},\n    @{Old = \"86\u00d780=6880\"; New = \"48\u00d718=864\"},\n    @{Old = \"88\u00d746=4048\"; New = \"23\u00d722=506\"},\n    @{Old = \"11\u00d753=583\"; New = \"82\u00d793=7626\"},\n    @{Old = \"94\u00d743=4042\"; New = \"20\u00d721=420\"},\n    @{Old = \"38\u00d715=570\"; New = \"84\u00d773=6132\"},\n    @{Old = \"71\u00d743=3053\"; New = \"45\u00d747=2115\"},\n    @{Old = \"46\u00d779=3634\"; New = \"76\u00d769=5244\"},\n    @{Old = \"33\u00d783=2739\"; New = \"33\u00d769=2277\"},\n    @{Old = \"55\u00d769=3795\"; New = \"13\u00d711=143\"}\n)\n\nforeach ($rep in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $rep.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $rep.New\n    $find.Execute(\n        $rep.Old,   # FindText\n        $false,     # MatchCase\n        $true,      # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $rep.New,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
